$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "model_11_8_0",
    "model_11_8_22",
    "model_11_8_21",
    "model_11_8_20",
    "model_11_8_19",
    "model_11_8_18",
    "model_11_8_17",
    "model_11_8_16",
    "model_11_8_15",
    "model_11_8_14",
    "model_11_8_13",
    "model_11_8_23",
    "model_11_8_12",
    "model_11_8_10",
    "model_11_8_9",
    "model_11_8_8",
    "model_11_8_7",
    "model_11_8_6",
    "model_11_8_5",
    "model_11_8_4",
    "model_11_8_3",
    "model_11_8_2",
    "model_11_8_1",
    "model_11_8_11",
    "model_11_8_24"
)

$values = @(0.3494677884409869, 0.2167805207458163, 0.221043372018899, 0.2227633730967914, 0.7199474573135376, 0.673530101776123, 0.9806371331214905, 0.8180509805679321)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
